# Auto-optimize exam scheduling: dynamically adjusts exams per slot (1-4)
# to guarantee all courses are scheduled within date range.
# Rewrites the Section_A and Section_B timetable grids (Mon-Fri columns
# B:F, class-slot rows 2,3,5,6,7,8) with the re-balanced schedule.

$wb = $excel.ActiveWorkbook

$wsA = $wb.Worksheets.Item("Section_A")
$wsB = $wb.Worksheets.Item("Section_B")

# --- Section_A -------------------------------------------------------
$wsA.Range("B2").Value = "Free"
$wsA.Range("C2").Value = "Free"
$wsA.Range("D2").Value = "CS307"
$wsA.Range("E2").Value = "CS307"
$wsA.Range("F2").Value = "EC303"

$wsA.Range("B3").Value = "EC262"
$wsA.Range("C3").Value = "EC264"
$wsA.Range("D3").Value = "Free"
$wsA.Range("E3").Value = "Free"
$wsA.Range("F3").Value = "Free"

$wsA.Range("B5").Value = "Free"
$wsA.Range("C5").Value = "EC304"
$wsA.Range("D5").Value = "EC303"
$wsA.Range("E5").Value = "EC264"
$wsA.Range("F5").Value = "EC262"

$wsA.Range("B6").Value = "EC304 (Tutorial)"
$wsA.Range("C6").Value = "Free"
$wsA.Range("D6").Value = "Free"
$wsA.Range("E6").Value = "Free"
$wsA.Range("F6").Value = "Free"

$wsA.Range("B7").Value = "CS307"
$wsA.Range("C7").Value = "EC303"
$wsA.Range("D7").Value = "EC262"
$wsA.Range("E7").Value = "EC304"
$wsA.Range("F7").Value = "EC264"

$wsA.Range("B8").Value = "EC303 (Tutorial)"
$wsA.Range("C8").Value = "Free"
$wsA.Range("D8").Value = "CS307 (Tutorial)"
$wsA.Range("E8").Value = "Free"
$wsA.Range("F8").Value = "Free"

# --- Section_B -------------------------------------------------------
$wsB.Range("B2").Value = "Free"
$wsB.Range("C2").Value = "EC303"
$wsB.Range("D2").Value = "Free"
$wsB.Range("E2").Value = "Free"
$wsB.Range("F2").Value = "EC264"

$wsB.Range("B3").Value = "CS307"
$wsB.Range("C3").Value = "EC264"
$wsB.Range("D3").Value = "CS307"
$wsB.Range("E3").Value = "CS307"
$wsB.Range("F3").Value = "EC262"

$wsB.Range("B5").Value = "Free"
$wsB.Range("C5").Value = "EC262"
$wsB.Range("D5").Value = "EC262"
$wsB.Range("E5").Value = "EC304"
$wsB.Range("F5").Value = "EC303"

$wsB.Range("B6").Value = "Free"
$wsB.Range("C6").Value = "Free"
$wsB.Range("D6").Value = "EC304 (Tutorial)"
$wsB.Range("E6").Value = "Free"
$wsB.Range("F6").Value = "Free"

$wsB.Range("B7").Value = "Free"
$wsB.Range("C7").Value = "Free"
$wsB.Range("D7").Value = "EC264"
$wsB.Range("E7").Value = "EC303"
$wsB.Range("F7").Value = "EC304"

$wsB.Range("B8").Value = "CS307 (Tutorial)"
$wsB.Range("C8").Value = "Free"
$wsB.Range("D8").Value = "Free"
$wsB.Range("E8").Value = "Free"
$wsB.Range("F8").Value = "EC303 (Tutorial)"
